$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Risk driver 1.1 answer options: change wording from "in total" to
# "on average per MT-member" for the four Management Team experience
# answers (rows 2-5, column E).
$ws.Range("E2").Value = "Management Team has on average per MT-member <1 years as a board member/leader/management team of innovative, circular and/or PaaS businesses"
$ws.Range("E3").Value = "Management Team has on average per MT-member 1-5 years of experience as a board member/leader/management team of innovative, circular and/or PaaS businesses"
$ws.Range("E4").Value = "Management Team has on average per MT-member 5-20 years of experience as a board member/leader/management team of innovative, circular and/or PaaS businesses"
$ws.Range("E5").Value = "Management Team has on average per MT-member >20 years of experience as a board member/leader/management team of innovative, circular and/or PaaS businesses"

# Update the active cell/selection to match the author's final cursor spot.
$ws.Range("H5").Select()
